$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1628.9546
$ws.Range("I80").Value = 1106.2142
$ws.Range("J80").Value = 2543.75
$ws.Range("K80").Value = 3318.6426
$ws.Range("L80").Value = 7631.25
$ws.Range("M80").Value = -2320.6426
$ws.Range("N80").Value = -9627.25
$ws.Range("H83").Value = 1628.9546
$ws.Range("I83").Value = 1106.2142
$ws.Range("J83").Value = 2543.75
$ws.Range("K83").Value = 9955.927799999999
$ws.Range("L83").Value = 22893.75
$ws.Range("M83").Value = -4963.927799999999
$ws.Range("N83").Value = -32877.75
$ws.Range("H98").Value = 2971332.5
$ws.Range("I98").Value = 3954530.2
$ws.Range("J98").Value = 915555.0600000001
$ws.Range("K98").Value = 3954530.2
$ws.Range("L98").Value = 915555.0600000001
$ws.Range("M98").Value = -3953032.2
$ws.Range("N98").Value = -918551.0600000001
$ws.Range("H106").Value = 501002.5
$ws.Range("I106").Value = 501002.5
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 501002.5
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -500371.5
$ws.Range("N106").Value = $null
$ws.Range("H122").Value = 2971332.5
$ws.Range("I122").Value = 3954530.2
$ws.Range("J122").Value = 915555.0600000001
$ws.Range("K122").Value = 11863590.6
$ws.Range("L122").Value = 2746665.18
$ws.Range("M122").Value = -11861140.6
$ws.Range("N122").Value = -2751565.18
$ws.Range("H137").Value = 1746275.5
$ws.Range("I137").Value = 3747806
$ws.Range("J137").Value = 3007.1936
$ws.Range("K137").Value = 11243418
$ws.Range("L137").Value = 9021.5808
$ws.Range("M137").Value = -11240868
$ws.Range("N137").Value = -14121.5808
$ws.Range("H138").Value = 1744.76
$ws.Range("I138").Value = 929.6
$ws.Range("J138").Value = 2094.1143
$ws.Range("K138").Value = 2788.8
$ws.Range("L138").Value = 6282.342900000001
$ws.Range("M138").Value = 2351.2
$ws.Range("N138").Value = -16562.3429

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6638.9043
$ws.Range("I32").Value = 5288.076
$ws.Range("K32").Value = 5288.076
$ws.Range("M32").Value = -5001.076
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080
$ws.Range("H132").Value = 2065.603
$ws.Range("I132").Value = 2144.1406
$ws.Range("J132").Value = 809
$ws.Range("K132").Value = 6432.4218
$ws.Range("L132").Value = 2427
$ws.Range("M132").Value = -3902.4218
$ws.Range("N132").Value = -7487

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2803.9707
$ws.Range("I134").Value = 2334.5715
$ws.Range("J134").Value = 4994.5
$ws.Range("K134").Value = 7003.7145
$ws.Range("L134").Value = 14983.5
$ws.Range("M134").Value = -4468.7145
$ws.Range("N134").Value = -20053.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 56155028
$ws.Range("I31").Value = 1113140.2
$ws.Range("J31").Value = 111196920
$ws.Range("K31").Value = 1113140.2
$ws.Range("L31").Value = 111196920
$ws.Range("M31").Value = -1112845.2
$ws.Range("N31").Value = -111197510
$ws.Range("H34").Value = 56155028
$ws.Range("I34").Value = 1113140.2
$ws.Range("J34").Value = 111196920
$ws.Range("K34").Value = 1113140.2
$ws.Range("L34").Value = 111196920
$ws.Range("M34").Value = -1112938.2
$ws.Range("N34").Value = -111197324
$ws.Range("H58").Value = 1955.2565
$ws.Range("I58").Value = 1935.4166
$ws.Range("J58").Value = 2193.3333
$ws.Range("K58").Value = 1935.4166
$ws.Range("L58").Value = 2193.3333
$ws.Range("M58").Value = -1732.4166
$ws.Range("N58").Value = -2599.3333
$ws.Range("H62").Value = 2597.5454
$ws.Range("I62").Value = 2481.8572
$ws.Range("J62").Value = 2800
$ws.Range("K62").Value = 2481.8572
$ws.Range("L62").Value = 2800
$ws.Range("M62").Value = -1857.8572
$ws.Range("N62").Value = -4048
$ws.Range("H65").Value = 2597.5454
$ws.Range("I65").Value = 2481.8572
$ws.Range("J65").Value = 2800
$ws.Range("K65").Value = 12409.286
$ws.Range("L65").Value = 14000
$ws.Range("M65").Value = -9289.286
$ws.Range("N65").Value = -20240
$ws.Range("H132").Value = 11685
$ws.Range("I132").Value = 2106.25
$ws.Range("J132").Value = 50000
$ws.Range("K132").Value = 6318.75
$ws.Range("L132").Value = 150000
$ws.Range("M132").Value = -3788.75
$ws.Range("N132").Value = -155060
$ws.Range("H134").Value = 3537
$ws.Range("I134").Value = 3526.6304
$ws.Range("J134").Value = 4014
$ws.Range("K134").Value = 10579.8912
$ws.Range("L134").Value = 12042
$ws.Range("M134").Value = -8044.8912
$ws.Range("N134").Value = -17112
$ws.Range("H136").Value = 1955.2565
$ws.Range("I136").Value = 1935.4166
$ws.Range("J136").Value = 2193.3333
$ws.Range("K136").Value = 5806.2498
$ws.Range("L136").Value = 6579.999899999999
$ws.Range("M136").Value = -3256.2498
$ws.Range("N136").Value = -11679.9999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3125113.5
$ws.Range("I2").Value = 345.5
$ws.Range("J2").Value = 3846213.8
$ws.Range("K2").Value = 2073
$ws.Range("L2").Value = 23077282.8
$ws.Range("M2").Value = -1960
$ws.Range("N2").Value = -23077508.8
$ws.Range("H131").Value = 6668270.5
$ws.Range("J131").Value = 1712.1692
$ws.Range("L131").Value = 5136.5076
$ws.Range("N131").Value = -15216.5076
$ws.Range("H140").Value = 5557203
$ws.Range("J140").Value = 1967.5358
$ws.Range("L140").Value = 5902.607400000001
$ws.Range("N140").Value = -16262.6074

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 21990.143
$ws.Range("I132").Value = 24232.598
$ws.Range("J132").Value = 4611.125
$ws.Range("K132").Value = 72697.79400000001
$ws.Range("L132").Value = 13833.375
$ws.Range("M132").Value = -70167.79400000001
$ws.Range("N132").Value = -18893.375

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1078.75
$ws.Range("I22").Value = 1328.3334
$ws.Range("J22").Value = 995.55554
$ws.Range("K22").Value = 1328.3334
$ws.Range("L22").Value = 995.55554
$ws.Range("M22").Value = -1033.3334
$ws.Range("N22").Value = -1585.55554
$ws.Range("H27").Value = 1078.75
$ws.Range("I27").Value = 1328.3334
$ws.Range("J27").Value = 995.55554
$ws.Range("K27").Value = 1328.3334
$ws.Range("L27").Value = 995.55554
$ws.Range("M27").Value = -1221.3334
$ws.Range("N27").Value = -1209.55554
$ws.Range("H40").Value = 8737.604499999999
$ws.Range("I40").Value = 9418.703
$ws.Range("K40").Value = 9418.703
$ws.Range("M40").Value = -9282.703
$ws.Range("H132").Value = 6468.533
$ws.Range("I132").Value = 5579.2
$ws.Range("K132").Value = 16737.6
$ws.Range("M132").Value = -14207.6
$ws.Range("H136").Value = 2988.4
$ws.Range("I136").Value = 1731.4255
$ws.Range("K136").Value = 5194.2765
$ws.Range("M136").Value = -2644.2765
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 21741974
$ws.Range("I132").Value = 26318346
$ws.Range("J132").Value = 4198.25
$ws.Range("K132").Value = 78955038
$ws.Range("L132").Value = 12594.75
$ws.Range("M132").Value = -78952508
$ws.Range("N132").Value = -17654.75
$ws.Range("H136").Value = 169910.44
$ws.Range("I136").Value = 230153.66
$ws.Range("J136").Value = 4241.5625
$ws.Range("K136").Value = 690460.98
$ws.Range("L136").Value = 12724.6875
$ws.Range("M136").Value = -687910.98
$ws.Range("N136").Value = -17824.6875
